# Update Active_Outages.xlsx - 6/18/2025, 4:33:45 PM
#
# Refreshes the "Elapsed Duration(Hrs)" figures across the regional tabs and
# consolidates the split R4 outage record on the R1 tab into a single row.

$wb = $excel.ActiveWorkbook

# --- Refresh elapsed-duration figures (Column G) on each region tab ---

$ws1 = $wb.Worksheets.Item("R1")
$ws1.Range("G2").Value = "3929:47:52"
$ws1.Range("G3").Value = "69:20:30"

$ws2 = $wb.Worksheets.Item("R2")
$ws2.Range("G2").Value = "12111:11:30"
$ws2.Range("G3").Value = "3240:54:59"
$ws2.Range("G4").Value = "479:06:33"

$ws4 = $wb.Worksheets.Item("R4")
$ws4.Range("G2").Value = "2957:01:19"
$ws4.Range("G3").Value = "184:13:34"
$ws4.Range("G4").Value = "72:25:59"
$ws4.Range("G5").Value = "70:03:32"

$ws5 = $wb.Worksheets.Item("R5")
$ws5.Range("G2").Value = "431:00:18"

$ws6 = $wb.Worksheets.Item("R6")
$ws6.Range("G2").Value = "71:32:36"

# --- Consolidate the R4 outage that was split across rows 4 & 5 of the R1 tab ---
# Row 5 carries the real Hub Site / Power Source / Battery Backup values for
# this outage; pull them into row 4, then drop the now-redundant row 5.

$ws1.Range("D4").Value = $ws1.Range("D5").Value
$ws1.Range("I4").Value = $ws1.Range("I5").Value
$ws1.Range("J4").Value = $ws1.Range("J5").Value

$ws1.Rows.Item(5).Delete()
